$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 51331.285
$ws.Range("I64").Value = 335533.34
$ws.Range("J64").Value = 3964.2778
$ws.Range("K64").Value = 335533.34
$ws.Range("L64").Value = 3964.2778
$ws.Range("M64").Value = -335285.34
$ws.Range("N64").Value = -4460.2778
# Row 67
$ws.Range("H67").Value = 51331.285
$ws.Range("I67").Value = 335533.34
$ws.Range("J67").Value = 3964.2778
$ws.Range("K67").Value = 335533.34
$ws.Range("L67").Value = 3964.2778
$ws.Range("M67").Value = -334675.34
$ws.Range("N67").Value = -5680.2778
# Row 86
$ws.Range("H86").Value = 4112.091
$ws.Range("I86").Value = 1387.7778
$ws.Range("K86").Value = 1387.7778
$ws.Range("M86").Value = -264.7778000000001
# Row 89
$ws.Range("H89").Value = 4112.091
$ws.Range("I89").Value = 1387.7778
$ws.Range("K89").Value = 6938.889
$ws.Range("M89").Value = -1322.889
# Row 98
$ws.Range("H98").Value = 8000
$ws.Range("I98").Value = 8000
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 8000
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -6502
$ws.Range("N98").ClearContents()
# Row 103
$ws.Range("H103").Value = 588.1
$ws.Range("I103").Value = 412.2
$ws.Range("J103").Value = 764
$ws.Range("K103").Value = 1236.6
$ws.Range("L103").Value = 2292
$ws.Range("M103").Value = -650.5999999999999
$ws.Range("N103").Value = -3464
# Row 122
$ws.Range("H122").Value = 8000
$ws.Range("I122").Value = 8000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 24000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -21550
$ws.Range("N122").ClearContents()
# Row 127
$ws.Range("H127").Value = 2302.8225
$ws.Range("I127").Value = 512.5
$ws.Range("J127").Value = 2426.2932
$ws.Range("K127").Value = 1537.5
$ws.Range("L127").Value = 7278.8796
$ws.Range("M127").Value = 3422.5
$ws.Range("N127").Value = -17198.8796

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 29
$ws.Range("H29").Value = 9266.666999999999
$ws.Range("J29").Value = 9266.666999999999
$ws.Range("L29").Value = 9266.666999999999
$ws.Range("N29").Value = -9882.666999999999
# Row 122
$ws.Range("H122").Value = 1315.2273
$ws.Range("I122").Value = 1252.75
$ws.Range("K122").Value = 3758.25
$ws.Range("M122").Value = -1308.25

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 17
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
# Row 99
$ws.Range("H99").Value = 17124.25
$ws.Range("I99").Value = 4160
$ws.Range("J99").Value = 24902.8
$ws.Range("K99").Value = 4160
$ws.Range("L99").Value = 24902.8
$ws.Range("M99").Value = -2662
$ws.Range("N99").Value = -27898.8
# Row 105
$ws.Range("H105").Value = 2218.8147
$ws.Range("I105").Value = 2303.9473
$ws.Range("J105").Value = 2016.625
$ws.Range("K105").Value = 2303.9473
$ws.Range("L105").Value = 2016.625
$ws.Range("M105").Value = -556.9472999999998
$ws.Range("N105").Value = -5510.625
# Row 126
$ws.Range("H126").Value = 17124.25
$ws.Range("I126").Value = 4160
$ws.Range("J126").Value = 24902.8
$ws.Range("K126").Value = 12480
$ws.Range("L126").Value = 74708.39999999999
$ws.Range("M126").Value = -10010
$ws.Range("N126").Value = -79648.39999999999

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 87
$ws.Range("H87").Value = 5687.375
$ws.Range("I87").Value = 4642.7144
$ws.Range("J87").Value = 13000
$ws.Range("K87").Value = 13928.1432
$ws.Range("L87").Value = 39000
$ws.Range("M87").Value = -12680.1432
$ws.Range("N87").Value = -41496
# Row 90
$ws.Range("H90").Value = 5687.375
$ws.Range("I90").Value = 4642.7144
$ws.Range("J90").Value = 13000
$ws.Range("K90").Value = 41784.4296
$ws.Range("L90").Value = 117000
$ws.Range("M90").Value = -35544.4296
$ws.Range("N90").Value = -129480
# Row 129
$ws.Range("H129").Value = 6225.4287
$ws.Range("J129").Value = 8522.666999999999
$ws.Range("L129").Value = 25568.001
$ws.Range("N129").Value = -35568.001

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 25
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("N25").ClearContents()
# Row 35
$ws.Range("H35").Value = 14319.333
$ws.Range("J35").Value = 14319.333
$ws.Range("L35").Value = 14319.333
$ws.Range("N35").Value = -14915.333
# Row 80
$ws.Range("H80").Value = 125126584
$ws.Range("I80").Value = 200201500
$ws.Range("J80").Value = 1735.3334
$ws.Range("K80").Value = 200201500
$ws.Range("L80").Value = 1735.3334
$ws.Range("M80").Value = -200200502
$ws.Range("N80").Value = -3731.3334
# Row 83
$ws.Range("H83").Value = 125126584
$ws.Range("I83").Value = 200201500
$ws.Range("J83").Value = 1735.3334
$ws.Range("K83").Value = 1001007500
$ws.Range("L83").Value = 8676.666999999999
$ws.Range("M83").Value = -1001002508
$ws.Range("N83").Value = -18660.667
# Row 122
$ws.Range("H122").Value = 1786.5264
$ws.Range("I122").Value = 2018.2142
$ws.Range("K122").Value = 6054.642599999999
$ws.Range("M122").Value = -3604.642599999999
# Row 126
$ws.Range("H126").Value = 4311.5625
$ws.Range("I126").Value = 4082.0833
$ws.Range("K126").Value = 12246.2499
$ws.Range("M126").Value = -9776.249899999999

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2071.8845
$ws.Range("I7").Value = 1404.0625
$ws.Range("K7").Value = 1404.0625
$ws.Range("M7").Value = -1292.0625
# Row 40
$ws.Range("H40").Value = 75140.36
$ws.Range("I40").Value = 147709.28
$ws.Range("J40").Value = 2571.4285
$ws.Range("K40").Value = 147709.28
$ws.Range("L40").Value = 2571.4285
$ws.Range("M40").Value = -147573.28
$ws.Range("N40").Value = -2843.4285
# Row 122
$ws.Range("H122").Value = 1617.3334
# Row 126
$ws.Range("H126").Value = 2071.8845
$ws.Range("I126").Value = 1404.0625
$ws.Range("K126").Value = 4212.1875
$ws.Range("M126").Value = -1742.1875
# Row 132
$ws.Range("H132").Value = 2842.2693
$ws.Range("I132").Value = 3051.6487
$ws.Range("J132").Value = 2325.8
$ws.Range("K132").Value = 9154.946100000001
$ws.Range("L132").Value = 6977.400000000001
$ws.Range("M132").Value = -6624.946100000001
$ws.Range("N132").Value = -12037.4

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 28
$ws.Range("H28").Value = 13450
$ws.Range("I28").Value = 7500
$ws.Range("J28").Value = 14937.5
$ws.Range("K28").Value = 7500
$ws.Range("L28").Value = 14937.5
$ws.Range("M28").Value = -7152
$ws.Range("N28").Value = -15633.5
# Row 29
$ws.Range("H29").Value = 200042370
$ws.Range("I29").Value = 500000900
$ws.Range("K29").Value = 500000900
$ws.Range("M29").Value = -500000610
# Row 122
$ws.Range("H122").Value = 1178.6
$ws.Range("I122").Value = 973.25
$ws.Range("K122").Value = 2919.75
$ws.Range("M122").Value = -469.75
# Row 126
$ws.Range("H126").Value = 1605.0714
$ws.Range("I126").Value = 1414.25
$ws.Range("J126").Value = 2750
$ws.Range("K126").Value = 4242.75
$ws.Range("L126").Value = 8250
$ws.Range("M126").Value = -1772.75
$ws.Range("N126").Value = -13190
